$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text (avoid Excel auto-numeric coercion),
# matching the original inlineStr text cells.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.071.10'
$ws.Range('E2').Value = '  -4.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.972.43'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.74'
$ws.Range('E5').Value = '  -3.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.84'
$ws.Range('E6').Value = '  +5.48%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +3.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.965.54'
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.130'
$ws.Range('E10').Value = '  -2.64%  '
$ws.Range('E11').Value = '  -5.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.13'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.462.12'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.84'
$ws.Range('E17').Value = '  +7.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.969.23'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '57.968.54'
$ws.Range('E19').Value = '  -4.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '422.34'
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.28'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.689'
$ws.Range('E22').Value = '  +3.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.02'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.10'
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.77'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.62'
$ws.Range('E29').Value = '  +3.85%  '
$ws.Range('E30').Value = '  +5.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.45'
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.11'
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0994'
$ws.Range('E33').Value = '  +6.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.68'
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('E35').Value = '  -1.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.948'
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0698'
$ws.Range('E37').Value = '  +5.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.69'
$ws.Range('E38').Value = '  -2.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.68'
$ws.Range('E39').Value = '  +3.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.61'
$ws.Range('E40').Value = '  +4.26%  '
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '378.94'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.706.23'
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.243'
$ws.Range('E46').Value = '  +2.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.56'
$ws.Range('E47').Value = '  +3.45%  '
$ws.Range('E48').Value = '  +2.64%  '
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.62'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('E51').Value = '  +0.07%  '
